# ---------------------------------------------------------------------------
# Rebuild the "error size / variation" side-table in columns L:P (rows 34-44)
# of Sheet1, matching the committed OOXML diff:
#  - drop the old title/footnote strings parked in column R (rows 29-32)
#  - relabel + recompute the variation matrix, splitting it into two
#    labeled blocks ("RMSE variation" / "MAE variation") each preceded by
#    a new "Error size" row and a merged section header
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PCT = "0.00%"

# ---------------------------------------------------------------------------
# 1. Remove the stray column-R notes (old shared strings 46, 56, 57, 58)
# ---------------------------------------------------------------------------
$ws.Range("R29").ClearContents()
$ws.Range("R30").ClearContents()
$ws.Range("R31").ClearContents()
$ws.Range("R32").ClearContents()

# ---------------------------------------------------------------------------
# 2. Fix up the SES/Croston/SBA/SBJ header row (row 30) - same text, the
#    underlying shared-string indices shift once the R-column strings go.
# ---------------------------------------------------------------------------
$ws.Range("M30").Value = "SES"
$ws.Range("N30").Value = "Croston"
$ws.Range("O30").Value = "SBA"
$ws.Range("P30").Value = "SBJ"

# ---------------------------------------------------------------------------
# 3. Break the old merges that spanned L34:L36 and L37:L39
# ---------------------------------------------------------------------------
$ws.Range("L34:L36").UnMerge()
$ws.Range("L37:L39").UnMerge()

# ---------------------------------------------------------------------------
# 4. Row 34 -> new "Error size" row: M..P = row33 - row32 (2-decimal style,
#    same look as the M31:P33 block above it)
# ---------------------------------------------------------------------------
$ws.Range("L34").Value = "Error size"
$ws.Range("L34").NumberFormat = "General"

$ws.Range("M34").Formula = "=M33-M32"
$ws.Range("N34").Formula = "=N33-N32"
$ws.Range("O34").Formula = "=O33-O32"
$ws.Range("P34").Formula = "=P33-P32"
$ws.Range("M34:P34").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 5. Row 35 -> merged section header "RMSE variation" (center / middle /
#    wrap, no number format) spanning L35:P35
# ---------------------------------------------------------------------------
$ws.Range("L35").Value = "RMSE variation"
$ws.Range("L35:P35").NumberFormat = "General"
$ws.Range("L35:P35").HorizontalAlignment = -4108
$ws.Range("L35:P35").VerticalAlignment = -4108
$ws.Range("L35:P35").WrapText = $true
$ws.Range("L35:P35").Merge()

# ---------------------------------------------------------------------------
# 6. Rows 36-39: RMSE variation matrix (based on row 33 RMSE values)
# ---------------------------------------------------------------------------
$ws.Range("L36").Value = "SES"
$ws.Range("M36").Value = ""
$ws.Range("N36").Formula = "=(M33-N33)/N33"
$ws.Range("O36").Formula = "=(M33-O33)/O33"
$ws.Range("P36").Formula = "=(M33-P33)/P33"

$ws.Range("L37").Value = "Croston"
$ws.Range("M37").Formula = "=(N33-M33)/M33"
$ws.Range("N37").Value = ""
$ws.Range("O37").Formula = "=(N33-O33)/O33"
$ws.Range("P37").Formula = "=(N33-P33)/P33"
$ws.Rows("37").RowHeight = 15

$ws.Range("L38").Value = "SBA"
$ws.Range("M38").Formula = "=(O33-M33)/M33"
$ws.Range("N38").Formula = "=(O33-N33)/N33"
$ws.Range("O38").Value = ""
$ws.Range("P38").Formula = "=(O33-P33)/P33"

$ws.Range("L39").Value = "SBJ"
$ws.Range("M39").Formula = "=(P33-M33)/M33"
$ws.Range("N39").Formula = "=(P33-N33)/N33"
$ws.Range("O39").Formula = "=(P33-O33)/O33"
$ws.Range("P39").Value = ""
$ws.Rows("39").RowHeight = 15

$ws.Range("L36:P39").NumberFormat = $PCT

# ---------------------------------------------------------------------------
# 7. Row 40 -> merged section header "MAE variation" spanning L40:P40
# ---------------------------------------------------------------------------
$ws.Range("L40").Value = "MAE variation"
$ws.Range("L40:P40").NumberFormat = $PCT
$ws.Range("L40:P40").HorizontalAlignment = -4108
$ws.Range("L40:P40").VerticalAlignment = -4108
$ws.Range("L40:P40").WrapText = $true
$ws.Range("L40:P40").Merge()

# ---------------------------------------------------------------------------
# 8. Rows 41-44: MAE variation matrix (based on row 32 MAE values)
# ---------------------------------------------------------------------------
$ws.Range("L41").Value = "SES"
$ws.Range("M41").Value = ""
$ws.Range("N41").Formula = "=(M32-N32)/N32"
$ws.Range("O41").Formula = "=(M32-O32)/O32"
$ws.Range("P41").Formula = "=(M32-P32)/P32"

$ws.Range("L42").Value = "Croston"
$ws.Range("M42").Formula = "=(N32-M32)/M32"
$ws.Range("N42").Value = ""
$ws.Range("O42").Formula = "=(N32-O32)/O32"
$ws.Range("P42").Formula = "=(N32-P32)/P32"

$ws.Range("L43").Value = "SBA"
$ws.Range("M43").Formula = "=(O32-M32)/M32"
$ws.Range("N43").Formula = "=(O32-N32)/N32"
$ws.Range("O43").Value = ""
$ws.Range("P43").Formula = "=(O32-P32)/P32"

$ws.Range("L44").Value = "SBJ"
$ws.Range("M44").Formula = "=(P32-M32)/M32"
$ws.Range("N44").Formula = "=(P32-N32)/N32"
$ws.Range("O44").Formula = "=(P32-O32)/O32"
$ws.Range("P44").Value = ""

$ws.Range("L41:P44").NumberFormat = $PCT

# ---------------------------------------------------------------------------
# 9. View state: selection on M34, scrolled so row 9 is at the top
# ---------------------------------------------------------------------------
$ws.Range("M34").Select()
$excel.ActiveWindow.ScrollRow = 9
